$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 131 first (higher index), then row 7, so indices don't shift under us.
$ws.Rows.Item(131).Delete()
$ws.Rows.Item(7).Delete()

# The old row 8 (now row 7) carried a yellow highlight on A7; clear it so no
# cell references the highlight fill anymore.
$ws.Range("A7").ClearFormats()

# Restore the view state (zoom, scroll position, selection) as last saved.
$excel.ActiveWindow.Zoom = 54
$ws.Range("M126").Select()
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 1
